$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Heading: "Objectives of the Study" -> "Objective of the Study"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Objectives of the Study", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Objective of the Study", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Replace the long "study aims to develop..." sentence
# ---------------------------------------------------------------------
$d.Content.Find.Execute("The study aims to develop a cloud-based learning center platform with mobile technology that helps learning centers facilitate and control their basic management processes.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "The study aims to develop a cloud-based learning center platform with mobile technology for administrative staff, educators, parents, and students.",
                         2) | Out-Null

# ---------------------------------------------------------------------
# 3-5. Append ";" to the end of the three objective bullet paragraphs
# ---------------------------------------------------------------------
$d.Content.Find.Execute("to gather data on the issues encountered by small and medium learning centers",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "to gather data on the issues encountered by small and medium learning centers;",
                         2) | Out-Null

$d.Content.Find.Execute("to design features on the app for both educators and learning centers",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "to design features on the app for both educators and learning centers;",
                         2) | Out-Null

$d.Content.Find.Execute("to define software requirements for both web and mobile development",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "to define software requirements for both web and mobile development;",
                         2) | Out-Null

# ---------------------------------------------------------------------
# 6. Replace the "Alleviating..." passage with the shorter replacement
# ---------------------------------------------------------------------
$apos = [char]0x2019
$newDelimText = "The app will regularly compare the applicant" + $apos + "s profile or details on every job hiring position and suggest the qualified applicants to the learning centers depending on the pre-set requirements and qualifications of the job hiring position."

$d.Content.Find.Execute("Alleviating the workload of Learning Centers is the primary aim of iLearnCentral. It will abridge the hiring and profiling of educators, scheduling, enrolment, etc. When hiring an educator, the app will suggest the best qualified applicant to the learning centers depending on the requirements and qualifications that were set.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         $newDelimText,
                         2) | Out-Null

# ---------------------------------------------------------------------
# 7. Insert "job seeking " before "educators can "
# ---------------------------------------------------------------------
$d.Content.Find.Execute("On the other hand, educators can apply for available learning center jobs through the app.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "On the other hand, job seeking educators can apply for available learning center jobs through the app.",
                         2) | Out-Null

Write-Host "done part 7"

# ---------------------------------------------------------------------
# 8. Append new sentence about Android to the end of that same paragraph
# ---------------------------------------------------------------------
$d.Content.Find.Execute("On the other hand, job seeking educators can apply for available learning center jobs through the app.  ",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "On the other hand, job seeking educators can apply for available learning center jobs through the app.  The app can run on Android 5.0 (Lollipop) or above.",
                         2) | Out-Null

Write-Host "done part 8"

# ---------------------------------------------------------------------
# 9. Replace the "Lastly, the app needs internet..." sentence
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Lastly, the app needs internet and Android 5.0 (Lollipop) or above to run. ",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Lastly, the app cannot fully function offline.",
                         2) | Out-Null

Write-Host "done part 9"
